$wb = $excel.ActiveWorkbook

# Locate the "study-data" worksheet (falls back to the first sheet if the
# name was already changed by a previous run).
$wsStudy = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "study-data") {
        $wsStudy = $sheet
        break
    }
}
if ($wsStudy -eq $null) {
    $wsStudy = $wb.Worksheets.Item(1)
}

# Rename it to "study_data"
$wsStudy.Name = "study_data"

# Make "study_data" the active sheet/tab (this clears tabSelected on every
# other sheet, e.g. "es_data"), and set its selection to E14
$wsStudy.Activate()
$wsStudy.Range("E14").Select()
